$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '59.662.51'
    'E2' = '  +3.62%  '
    'D3' = '3.018.57'
    'E3' = '  +3.17%  '
    'E4' = '  -0.01%  '
    'D5' = '565.40'
    'E5' = '  +3.17%  '
    'D6' = '140.06'
    'E6' = '  +7.85%  '
    'E7' = '  -0.02%  '
    'D8' = '0.521'
    'E8' = '  +2.15%  '
    'D9' = '3.010.79'
    'E9' = '  +3.12%  '
    'E10' = '  +6.24%  '
    'D11' = '5.29'
    'E11' = '  +11.33%  '
    'D12' = '0.462'
    'E12' = '  +3.56%  '
    'E13' = '  +5.72%  '
    'D14' = '34.04'
    'E14' = '  +3.83%  '
    'E15' = '  +1.85%  '
    'D16' = '3.514.67'
    'E16' = '  +3.37%  '
    'D17' = '7.27'
    'E17' = '  +6.46%  '
    'D18' = '3.017.82'
    'E18' = '  +3.32%  '
    'D19' = '59.644.95'
    'E19' = '  +3.60%  '
    'D20' = '436.62'
    'E20' = '  +4.71%  '
    'D21' = '13.70'
    'E22' = '  +6.32%  '
    'D23' = '7.14'
    'E23' = '  +2.78%  '
    'E24' = '  +2.14%  '
    'D25' = '80.86'
    'E25' = '  +1.40%  '
    'E26' = '  +0.10%  '
    'D27' = '2.26'
    'E27' = '  +14.07%  '
    'E28' = '  +0.02%  '
    'E29' = '  +3.64%  '
    'D30' = '7.85'
    'E30' = '  +5.76%  '
    'D31' = '26.05'
    'E31' = '  +3.56%  '
    'D32' = '6.29'
    'E32' = '  +5.99%  '
    'E33' = '  +5.26%  '
    'D34' = '0.0₃0788'
    'E34' = '  +16.47%  '
    'D35' = '1.01'
    'E35' = '  +7.19%  '
    'D36' = '5.94'
    'E36' = '  +5.42%  '
    'E37' = '  +2.69%  '
    'E38' = '  +2.80%  '
    'D39' = '8.67'
    'E39' = '  -0.10%  '
    'E40' = '  +9.12%  '
    'D41' = '404.53'
    'E41' = '  +8.12%  '
    'D42' = '0.0356'
    'E42' = '  +3.55%  '
    'D43' = '2.784.53'
    'E43' = '  +4.71%  '
    'E44' = '  +0.45%  '
    'E45' = '  +7.01%  '
    'E46' = '  +0.00%  '
    'D47' = '123.70'
    'E47' = '  +1.30%  '
    'B48' = 'Fetch.AI'
    'C48' = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    'D48' = '2.03'
    'E48' = '  +3.43%  '
    'B49' = 'Stellar'
    'C49' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D49' = '0.111'
    'E49' = '  +1.81%  '
    'D50' = '33.89'
    'E50' = '  +2.20%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
